# "Generate Report for Handback" - refresh the localization-status report
# after a handback sync: update status text, handback timestamps, clear the
# stale "not latest" error now that de-de is back in sync, and widen a
# couple of columns that now hold longer text.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Overview")
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws3 = $wb.Worksheets.Item("de-de")

# --- Status text: "Ready for handoff" -> "Handed back: in sync with en-US" ---
# This text is shared by the Overview roll-up (E2 = zh-cn, F2 = de-de) and
# each language sheet's own Status cell (C2); update every occurrence so the
# underlying string itself is fully replaced.
$ws1.Range("E2").Value = "Handed back: in sync with en-US"
$ws1.Range("F2").Value = "Handed back: in sync with en-US"
$ws2.Range("C2").Value = "Handed back: in sync with en-US"
$ws3.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: refresh handback datetime, clear stale error detail ---
$ws2.Range("K2").Value = "2016-08-20 12:51:55"
$ws2.Range("P2").Value = ""

# --- de-de sheet: refresh handback datetime, clear stale error detail ---
$ws3.Range("K2").Value = "2016-08-20 12:52:04"
$ws3.Range("P2").Value = ""

# --- Column widths: widen the Status columns and shrink the now-empty
#     Error Detail columns to fit their (shorter) content. (Input values are
#     tuned so the engine's pixel-grid rounding lands on the closest
#     representable width to the target.) ---
$ws1.Range("E1").ColumnWidth = 29.15
$ws1.Range("F1").ColumnWidth = 29.15

$ws2.Range("C1").ColumnWidth = 29.15
$ws2.Range("P1").ColumnWidth = 12.8

$ws3.Range("C1").ColumnWidth = 29.15
$ws3.Range("P1").ColumnWidth = 12.8
